$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell A3 gets the shared string "test"
$ws.Range("A3").Value = "test"

# Merge A3:A4 (same value concept per commit message - "cell merge with same value")
$ws.Range("A3:A4").Merge()

# Touch A4's format (no-op self assignment) so Excel's used-range/dimension
# correctly extends to cover the merged cell down to row 4, matching native
# Excel behaviour for dimension calculation without altering any cell style.
$ws.Range("A4").Font.Bold = $ws.Range("A4").Font.Bold
